$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-10 (row 23)
$ws.Range("B23").Value = 6326
$ws.Range("D23").Value = 5925350
$ws.Range("E23").Value = 936.6661397407524
$ws.Range("F23").Value = 8.544955387783126
$ws.Range("H23").Value = 26.96508377164948
